# Integration of "GT Structure" content: update metadata (Date/Description)
# and replace the Organization type concept list on the "Include #0" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("Metadata"): update Date and Description values -----------
$wsMeta = $wb.Worksheets.Item(1)
$wsMeta.Cells.Item(8, 2).Value = "2026-01-28T10:29:57+00:00"
$wsMeta.Cells.Item(13, 2).Value = "Types de structures"

# --- Sheet 2 ("Include #0"): update the concept / description list ------
$wsInc = $wb.Worksheets.Item(2)

# Rows 3-11 currently hold:
#   3  LEGAL-ENTITY       / Entité juridique
#   4  GEOGRAPHICAL-ENTITY/ Entité géographique
#   5  GROUP              / Groupe privé / hospitalier
#   6  HEBERGEMENT        / Hébergement
#   7  SOIN               / Soin
#   8  ADMINISTRATIF      / Administration
#   9  MEDICAL            / Médical
#   10 TECHNIQUE          / Technique
#   11 MEDICOTEC          / medico technique (radio, scanner …)
# Overwrite them in place (reusing existing row styles) with the new set of
# concepts: GROUP, STRUCT-INTERNE, SECTEUR, DEPARTEMENT, SERVICE, UM, POLE,
# CENTRE-RESP, CENTRE-ACTIVITE.
$wsInc.Cells.Item(3, 1).Value = "GROUP"
$wsInc.Cells.Item(3, 2).Value = "Groupe privé / hospitalier"
$wsInc.Cells.Item(4, 1).Value = "STRUCT-INTERNE"
$wsInc.Cells.Item(4, 2).Value = "Structure interne"
$wsInc.Cells.Item(5, 1).Value = "SECTEUR"
$wsInc.Cells.Item(5, 2).Value = "Secteur"
$wsInc.Cells.Item(6, 1).Value = "DEPARTEMENT"
$wsInc.Cells.Item(6, 2).Value = "Département"
$wsInc.Cells.Item(7, 1).Value = "SERVICE"
$wsInc.Cells.Item(7, 2).Value = "Service"
$wsInc.Cells.Item(8, 1).Value = "UM"
$wsInc.Cells.Item(8, 2).Value = "Unité médicale"
$wsInc.Cells.Item(9, 1).Value = "POLE"
$wsInc.Cells.Item(9, 2).Value = "Pole"
$wsInc.Cells.Item(10, 1).Value = "CENTRE-RESP"
$wsInc.Cells.Item(10, 2).Value = "Centre de responsabilité"
$wsInc.Cells.Item(11, 1).Value = "CENTRE-ACTIVITE"
$wsInc.Cells.Item(11, 2).Value = "Centre d'activité"

# The rows that followed (old SERVICE, UM, UAC, POLE, CENTRE-RESP,
# CENTRE-ACTIVITE - 7 rows starting at row 12 now that UF/etc. content was
# overwritten above) are no longer needed: their content has already been
# relocated above, so delete the 7 now-redundant rows. This shifts the
# trailing empty row and "System URI" row up without touching their styles.
for ($i = 0; $i -lt 7; $i++) {
    $wsInc.Rows.Item(12).Delete()
}
